$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q1" sheet right before "总计", by copying the
#    "2021-Q4" sheet (same 8-column fund-holding layout/styles) and then
#    trimming / overwriting its data.
# ---------------------------------------------------------------------------
$src   = $wb.Worksheets.Item("2021-Q4")
$total = $wb.Worksheets.Item("总计")
$src.Copy($total)
$newSheet = $wb.ActiveSheet
$newSheet.Name = "2022-Q1"

# Re-fetch "总计" by name: the $total reference above gets reseated to the
# freshly-inserted sheet once Copy() shifts sheet positions around it.
$total = $wb.Worksheets.Item("总计")

# The source sheet (2021-Q4) has 41 data rows; the new quarter only needs 15
# (rows 2-16), so drop the extra rows entirely (keeps dimension A1:H16).
$newSheet.Rows("17:41").Delete()

# Fund-holding rows for 2022-Q1 (B:基金代码 C:基金名称 D:基金规模
# E:股票总仓位 F:仓位占比 G:持有市值(亿元) H:仓位排名). Column A (0-based
# index) already holds the correct 0..14 sequence copied from 2021-Q4.
$rows = @(
    @("290011", "泰信中小盘精选混合",         "13.68", "94.64", "8.99", "1.2298", 6),
    @("610002", "信达澳银精华配置混合",        "18.32", "71.31", "2.55", "0.4672", 2),
    @("610001", "信达澳银领先增长混合",        "8.99",  "94.02", "4.48", "0.4028", 3),
    @("013840", "银华集成电路混合A",          "8.32",  "71.47", "3.53", "0.2937", 7),
    @("001970", "泰信鑫选灵活配置混合A",       "3.03",  "93.96", "8.70", "0.2636", 5),
    @("002580", "泰信鑫选灵活配置混合C",       "2.04",  "93.96", "8.70", "0.1775", 5),
    @("001125", "博时互联网主题灵活配置混合",   "5.34",  "84.75", "2.54", "0.1356", 9),
    @("506008", "长城科创两年定开混合A",       "3.57",  "62.09", "2.20", "0.0785", 8),
    @("006813", "博时汇悦回报混合",            "1.11",  "84.25", "3.08", "0.0342", 8),
    @("013841", "银华集成电路混合C",          "0.72",  "71.47", "3.53", "0.0254", 7),
    @("003659", "山西证券策略精选灵活配置混合", "0.31",  "84.52", "4.11", "0.0127", 3),
    @("004930", "华润元大价值优选混合A",       "0.32",  "65.19", "3.26", "0.0104", 10),
    @("004931", "华润元大价值优选混合C",       "0.18",  "65.19", "3.26", "0.0059", 10),
    @("012793", "长城科创两年定开混合C",       "0.12",  "62.09", "2.20", "0.0026", 8),
    @("002194", "北信瑞丰稳定增强偏债混合",     "0.04",  "22.00", "2.26", "0.0009", 2)
)

$r = 2
foreach ($row in $rows) {
    # Force B:G to stay text (matches source file, and protects fund codes
    # like "013840" / decimal-looking strings like "13.68" from being
    # silently re-typed as numbers).
    $newSheet.Range("B$r`:G$r").NumberFormat = "@"
    $newSheet.Range("B$r").Value = $row[0]
    $newSheet.Range("C$r").Value = $row[1]
    $newSheet.Range("D$r").Value = $row[2]
    $newSheet.Range("E$r").Value = $row[3]
    $newSheet.Range("F$r").Value = $row[4]
    $newSheet.Range("G$r").Value = $row[5]
    $newSheet.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing quarters down by one row and renumbering the index column.
# ---------------------------------------------------------------------------
$total.Range("A6").Copy()
$total.Range("A7").PasteSpecial(-4122) | Out-Null   # xlPasteFormats (extend the styled A column by one row)

$summary = @(
    @("2022-Q1", 15, 3.14),
    @("2021-Q4", 40, 6.25),
    @("2021-Q3", 12, 4.33),
    @("2021-Q2", 11, 8.21),
    @("2021-Q1", 17, 10.65),
    @("2020-Q4", 7, 3.87)
)

$r = 2
$idx = 0
foreach ($row in $summary) {
    $total.Range("A$r").Value = $idx
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
    $r = $r + 1
    $idx = $idx + 1
}
